# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect the newly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    2  = 1147
    3  = 414
    4  = 258
    6  = 5
    7  = 12162
    8  = 52
    9  = 9
    10 = 114
    11 = 11942
    12 = 4784
    13 = 871
    14 = 103
    15 = 40
    18 = 939
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# --- Sheet "全部类型" (all categories) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    2  = 1147
    3  = 414
    4  = 258
    8  = 5
    9  = 12162
    10 = 52
    11 = 9
    12 = 114
    13 = 11942
    14 = 4784
    15 = 873
    16 = 103
    17 = 40
    20 = 939
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
